$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.838.03'
$ws.Range('E2').Value = '  +1.08%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.666.50'
$ws.Range('E3').Value = '  +2.39%  '

$ws.Range('E4').Value = '  +0.42%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.71'
$ws.Range('E5').Value = '  +0.99%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.529'
$ws.Range('E6').Value = '  +5.55%  '

$ws.Range('E7').Value = '  +0.30%  '

$ws.Range('E8').Value = '  +2.65%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0619'
$ws.Range('E9').Value = '  +1.58%  '

$ws.Range('E10').Value = '  +4.89%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0891'
$ws.Range('E11').Value = '  +4.05%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.904.44'
$ws.Range('E12').Value = '  +2.54%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.697.52'
$ws.Range('E13').Value = '  +4.31%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.08'
$ws.Range('E14').Value = '  +0.84%  '

$ws.Range('E15').Value = '  +1.69%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.56'
$ws.Range('E16').Value = '  +2.59%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.887.22'
$ws.Range('E17').Value = '  +1.19%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '231.60'
$ws.Range('E18').Value = '  -1.04%  '

$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.78'
$ws.Range('E19').Value = '  -0.40%  '

$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0734'
$ws.Range('E20').Value = '  +1.30%  '

$ws.Range('E21').Value = '  +0.40%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.44'
$ws.Range('E22').Value = '  +2.65%  '

$ws.Range('B23').Value = 'Avalanche'
$ws.Range('C23').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.19'
$ws.Range('E23').Value = '  +0.37%  '

$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.21'
$ws.Range('E24').Value = '  +0.42%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.72'
$ws.Range('E25').Value = '  -0.18%  '

$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.116'
$ws.Range('E26').Value = '  +2.80%  '

$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.13'
$ws.Range('E27').Value = '  +0.78%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.86'
$ws.Range('E28').Value = '  +1.13%  '

$ws.Range('E29').Value = '  +0.23%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0496'
$ws.Range('E30').Value = '  +0.25%  '

$ws.Range('E31').Value = '  +1.21%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.31'
$ws.Range('E32').Value = '  +1.72%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.455.09'
$ws.Range('E33').Value = '  -4.62%  '

$ws.Range('E34').Value = '  +4.47%  '

$ws.Range('E35').Value = '  +5.13%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.43'
$ws.Range('E36').Value = '  +0.36%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.897'
$ws.Range('E37').Value = '  +7.21%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.565'
$ws.Range('E38').Value = '  -0.76%  '

$ws.Range('E39').Value = '  +0.98%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.03'
$ws.Range('E40').Value = '  +2.80%  '

$ws.Range('E41').Value = '  +0.31%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.29'
$ws.Range('E42').Value = '  +3.58%  '

$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.970'
$ws.Range('E43').Value = '  +6.79%  '

$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '65.59'
$ws.Range('E44').Value = '  +3.71%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.810.83'
$ws.Range('E45').Value = '  +2.38%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.777'
$ws.Range('E46').Value = '  +2.11%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.65'
$ws.Range('E47').Value = '  +0.99%  '

$ws.Range('E48').Value = '  +0.63%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0999'
$ws.Range('E49').Value = '  +3.46%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0507'
$ws.Range('E50').Value = '  +1.04%  '

$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₇0973'
$ws.Range('E51').Value = '  -6.71%  '
